# Commit: "added colour blindness test"
#
# A new slide (Ishihara colour-blindness plates instructions) is inserted
# right after the current first slide. The new slide carries the text that
# used to live on slide 1 ("grey line drawings ..."), and slide 1 itself is
# rewritten in place with the new colour-blindness instructions. Every other
# slide keeps its content/order unchanged, just shifted down by one position.

$p = $ppt.ActivePresentation

# --- 1. Duplicate slide 1. PowerPoint inserts the duplicate immediately
#        after the original, so it becomes the new slide #2 and keeps the
#        original "grey line drawings" wording/formatting/position intact.
$slide1 = $p.Slides.Item(1)
$slide1.Duplicate() | Out-Null

# --- 2. Rewrite the (original) slide 1 textbox with the new colour
#        blindness instructions, and resize/move it to its new, smaller
#        placement on the slide.
$shape = $p.Slides.Item(1).Shapes.Item(1)
$shape.TextFrame.TextRange.Text = "You will be shown 15 plates. Please enter the number you see on each plate followed by Enter key.`rPress the spacebar to begin"

$shape.Left = 381817 / 12700
$shape.Top = 2828836 / 12700
$shape.Width = 8380366 / 12700
$shape.Height = 1200329 / 12700
